$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.962721347808838
$ws.Range("B1").Value = 3.66393518447876
$ws.Range("C1").Value = 2.893548488616943
$ws.Range("D1").Value = 2.350312948226929
$ws.Range("E1").Value = 1.487265944480896
